$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.188.65"
$ws.Range("E2").Value = "  +5.31%  "

$ws.Range("D3").Value = "2.243.39"
$ws.Range("E3").Value = "  +4.78%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.24"
$ws.Range("E5").Value = "  +6.72%  "

$ws.Range("E6").Value = "  +2.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.22"
$ws.Range("E7").Value = "  +9.70%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +7.18%  "

$ws.Range("E10").Value = "  +8.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  +4.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.93"
$ws.Range("E12").Value = "  +5.65%  "

$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").Value = "2.580.73"
$ws.Range("E14").Value = "  +5.03%  "

$ws.Range("E15").Value = "  +2.18%  "

$ws.Range("D16").Value = "2.240.39"
$ws.Range("E16").Value = "  +5.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +2.54%  "

$ws.Range("D18").Value = "43.065.17"
$ws.Range("E18").Value = "  +5.48%  "

$ws.Range("E19").Value = "  +7.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.25"
$ws.Range("E20").Value = "  +3.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +5.53%  "

$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").Value = "  +18.79%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "229.84"
$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.80"
$ws.Range("E26").Value = "  +3.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.44"
$ws.Range("E27").Value = "  +4.89%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.36"
$ws.Range("E28").Value = "  +29.36%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +5.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +4.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.96"
$ws.Range("E31").Value = "  +2.52%  "

$ws.Range("E32").Value = "  +3.87%  "

$ws.Range("E33").Value = "  +8.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.26"
$ws.Range("E34").Value = "  +5.14%  "

$ws.Range("E35").Value = "  +2.76%  "

$ws.Range("E36").Value = "  +10.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.49"
$ws.Range("E37").Value = "  +11.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0333"
$ws.Range("E38").Value = "  +20.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.19"
$ws.Range("E39").Value = "  +14.57%  "

$ws.Range("E40").Value = "  +4.77%  "

$ws.Range("E41").Value = "  +12.37%  "

$ws.Range("E42").Value = "  +4.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.49"
$ws.Range("E43").Value = "  +4.16%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.68"
$ws.Range("E44").Value = "  +7.04%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.29"
$ws.Range("E45").Value = "  +9.01%  "

$ws.Range("E46").Value = "  +35.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0993"
$ws.Range("E47").Value = "  +5.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +14.08%  "

$ws.Range("E49").Value = "  +4.72%  "

$ws.Range("E50").Value = "  +5.53%  "

$ws.Range("E51").Value = "  +3.83%  "
